$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
# A2/B2 already carry the "Text" style (s=4); just change their values.
$ws.Range("A2").Value = "(unique value filled in by the test)"
$ws.Range("B2").Value = "(unique value filled in by the test)"

# C2/D2 need to pick up the "Text" number format (same style class as col A/B, numFmt 49)
# and new values.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "(unique value filled in by the test)"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "(same as broad sample)"

# O2 value change only (rest of the row is untouched).
$ws.Range("O2").Value = "lsid:1"

# --- Row 3: brand-new data row, replacing the old placeholder/formatting-only row ---
# Clear the stale height + the one styled (wrap-text) cell first.
$ws.Rows.Item(3).AutoFit()
$ws.Range("I3").Style = "Normal"

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "(unique value filled in by the test)"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "(unique value filled in by the test)"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "(unique value filled in by the test)"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "(unique value filled in by the test)"

$ws.Range("E3").Value = "Illumina_P5-Nijow_P7-Waren"
$ws.Range("H3").Value = "DEV-6796"
$ws.Range("I3").Value = "DEV-6815, DEV-6816"
$ws.Range("J3").Value = "COLB-124"
$ws.Range("K3").Value = "COLAB-P-235"
$ws.Range("L3").Value = "BP-ID-568"
$ws.Range("M3").Value = "F"
$ws.Range("N3").Value = "Feline"
$ws.Range("O3").Value = "lsid:2"
$ws.Range("P3").Value = 62
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 4

# --- Column widths (best-fit recalculated by Excel after the content changed) ---
$ws.Columns.Item(1).ColumnWidth = 10.830729166666666
$ws.Columns.Item(2).ColumnWidth = 20.998697916666668
$ws.Columns.Item(3).ColumnWidth = 12.998697916666666
$ws.Columns.Item(4).ColumnWidth = 12.166666666666666
$ws.Columns.Item(5).ColumnWidth = 22.498697916666668
$ws.Columns.Item(6).ColumnWidth = 24.830729166666668
$ws.Columns.Item(7).ColumnWidth = 24.830729166666668
$ws.Columns.Item(8).ColumnWidth = 9.166666666666666
$ws.Columns.Item(9).ColumnWidth = 15.998697916666666
$ws.Columns.Item(10).ColumnWidth = 17.998697916666668
$ws.Columns.Item(11).ColumnWidth = 20.830729166666668
$ws.Columns.Item(12).ColumnWidth = 15.830729166666666
$ws.Columns.Item(14).ColumnWidth = 7.666666666666667
$ws.Columns.Item(15).ColumnWidth = 5.330729166666667
$ws.Columns.Item(16).ColumnWidth = 7.666666666666667
$ws.Columns.Item(17).ColumnWidth = 14.166666666666666
$ws.Columns.Item(18).ColumnWidth = 12.498697916666666

# --- Selection moves to G3 ---
$ws.Range("G3").Select()
